$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.332649053379625
$ws.Range("D2").Value = 0.1962859687175309

$ws.Range("C3").Value = 1.05205638241968
$ws.Range("D3").Value = 0.3041927058056277

$ws.Range("C4").Value = 0.1078463366202454
$ws.Range("D4").Value = 0.9150951092465252

$ws.Range("C5").Value = 2.825770296413805
$ws.Range("D5").Value = 0.009841598548882224

$ws.Range("C6").Value = -0.2014396957651644
$ws.Range("D6").Value = 0.8422056787874563

$ws.Range("C7").Value = -1.102182581385788
$ws.Range("D7").Value = 0.282294425541433

$ws.Range("C8").Value = 1.613316690742703
$ws.Range("D8").Value = 0.1209292460206937

$ws.Range("C9").Value = -1.102146646766532
$ws.Range("D9").Value = 0.2823097043932687

$ws.Range("C10").Value = 1.444050001177935
$ws.Range("D10").Value = 0.162816200751013

$ws.Range("C11").Value = 1.756394177002759
$ws.Range("D11").Value = 0.09293393441114683
